# Auto-generated: refresh the cryptos table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.409.21'
$ws.Range("E2").Value = '  +3.79%  '

$ws.Range("D3").Value = '3.487.72'
$ws.Range("E3").Value = '  +2.69%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.29'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.45'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.52%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.609'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +12.64%  '

$ws.Range("D9").Value = '3.489.45'
$ws.Range("E9").Value = '  +2.72%  '

$ws.Range("E10").Value = '  -0.89%  '

$ws.Range("E11").Value = '  +3.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.449'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.05%  '

$ws.Range("D13").Value = '4.092.92'
$ws.Range("E13").Value = '  +2.66%  '

$ws.Range("E14").Value = '  +1.05%  '

$ws.Range("E15").Value = '  +3.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.69'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.76%  '

$ws.Range("D17").Value = '65.439.57'
$ws.Range("E17").Value = '  +3.64%  '

$ws.Range("D18").Value = '3.510.16'
$ws.Range("E18").Value = '  +2.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.48'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '387.40'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.29'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.68%  '

$ws.Range("E23").Value = '  +4.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.25'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.995'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("E26").Value = '  +5.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.98%  '

$ws.Range("E28").Value = '  +2.34%  '

$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  +11.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.27'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.68'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.05%  '

$ws.Range("E34").Value = '  +8.44%  '

$ws.Range("E35").Value = '  +8.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.77'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.26%  '

$ws.Range("E37").Value = '  +6.19%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0779'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.78%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '2.991.06'
$ws.Range("E39").Value = '  +3.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.26'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.25%  '

$ws.Range("E41").Value = '  +3.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.59'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.17%  '

$ws.Range("E43").Value = '  +2.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.65'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.98%  '

$ws.Range("E45").Value = '  +3.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.45'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +9.09%  '

$ws.Range("E47").Value = '  +4.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '322.21'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +12.21%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.76'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.18%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.110'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +7.32%  '

$ws.Range("E51").Value = '  +3.13%  '
